{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Implements three textual revisions to the document (matching the\n// supplied OOXML diff):\n//   1. \"DbContext\" paragraph \u2014 re-typing of \"class'dan\" (text unchanged).\n//   2. \"Configuration\" paragraph \u2014 re-typing of \"class'dan\" (text unchanged).\n//   3. \"AutoMapper\" section \u2014 expanded / rewritten explanation split\n//      across several paragraphs (new paragraphs inserted).\n\nconst RTL_APOSTROPHE = \"\\u2019\"; // \u2019 (RIGHT SINGLE QUOTATION MARK) used throughout the doc\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Helper: re-type a known substring inside a paragraph's range so the\n// run carrying it is (re)written \u2014 used for the two small \"class'dan\"\n// edits where the final text is unchanged but the run was re-entered.\n// ---------------------------------------------------------------------\nasync function retypeSubstring(paragraphText_contains, needle) {\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    const r = results.items[i];\n    const paraRange = r.paragraphs.getFirstOrNullObject();\n    paraRange.load(\"text\");\n    await context.sync();\n    if (!paraRange.isNullObject && paraRange.text.indexOf(paragraphText_contains) !== -1) {\n      r.insertText(needle, Word.InsertLocation.replace);\n      await context.sync();\n      return true;\n    }\n  }\n  return false;\n}\n\n// 1) \"\u0130lk olaraq DbContext class'dan miras alan ...\" paragraph.\nawait retypeSubstring(\"\u0130lk olaraq DbContext\", \"class\" + RTL_APOSTROPHE + \"dan\");\n\n// 2) \".cs file daxilind\u0259 ... Configuration class'dan istifad\u0259 etm\u0259kdir.\" paragraph.\nawait retypeSubstring(\"Configuration\", \"class\" + RTL_APOSTROPHE + \"dan\");\n\n// ---------------------------------------------------------------------\n// 3) AutoMapper section rewrite.\n// ---------------------------------------------------------------------\nconst paras2 = body.paragraphs;\nparas2.load(\"text\");\nawait context.sync();\n\nlet idxIntro = -1; // \"Biz project daxilind\u0259 Front'a Entity'l\u0259ri yox Model ...\"\nlet idxBody = -1; // \"AutoMapper v\u0259 AutoMapper.Extensions.Microsoft.DependencyInjection ...\"\nfor (let i = 0; i < paras2.items.length; i++) {\n  const t = paras2.items[i].text;\n  if (idxIntro === -1 && t.indexOf(\"Biz project daxilind\u0259 Front\") !== -1 && t.indexOf(\"Entity\") !== -1) {\n    idxIntro = i;\n  }\n  if (idxBody === -1 && t.indexOf(\"AutoMapper v\u0259 AutoMapper.Extensions.Microsoft.DependencyInjection\") !== -1) {\n    idxBody = i;\n  }\n}\n\nif (idxIntro === -1 || idxBody === -1) {\n  throw new Error(\"Could not locate AutoMapper paragraphs (idxIntro=\" + idxIntro + \", idxBody=\" + idxBody + \")\");\n}\n\nconst paraIntro = paras2.items[idxIntro]; // becomes target[0]\nconst paraBlank1 = paras2.items[idxIntro + 1]; // already empty, stays target[1]\nconst paraBody = paras2.items[idxBody]; // becomes target[7] (reused)\n\nconst newTexts = [\n  \"Biz security c\u0259h\u0259td\u0259n project daxilind\u0259 Front\" + RTL_APOSTROPHE + \"a database\" + RTL_APOSTROPHE + \"d\u0259n g\u0259l\u0259n Entity object\" + RTL_APOSTROPHE + \"l\u0259ri bir ba\u015fa olaraq g\u00f6nd\u0259rm\u0259yimiz t\u0259hl\u00fck\u0259lidir. \u00c7\u00fcnki biz bu object\" + RTL_APOSTROPHE + \"l\u0259ri bir ba\u015fa \u00f6t\u00fcrs\u0259k \u0259l\u00e7atan olmamal\u0131 bir s\u0131ra datalar\u0131 da (property\" + RTL_APOSTROPHE + \"l\u0259ri) g\u00f6nd\u0259rmi\u015f olar\u0131q. Bunun qar\u015f\u0131s\u0131n\u0131 almaq \u00fc\u00e7\u00fcn ViewModel v\u0259 ya DTO class\" + RTL_APOSTROPHE + \"lar\u0131 yarad\u0131b m\u00fc\u0259yy\u0259n olunmu\u015f property\" + RTL_APOSTROPHE + \"l\u0259ri saxlay\u0131r\u0131q.\",\n  \"H\u0259r d\u0259f\u0259 datan\u0131n \u00f6t\u00fcr\u00fclm\u0259si v\u0259 ya q\u0259bul olunmas\u0131 zaman\u0131 yeni ViewModel & DTO v\u0259 ya Entity object yarad\u0131lmas\u0131 v\u0259 d\u0259y\u0259rl\u0259rin m\u0259nims\u0259dilm\u0259si h\u0259m kod t\u0259krar\u0131na g\u0259tirib \u00e7\u0131xar\u0131r h\u0259m d\u0259 vaxt itkisin\u0259. Bu problemin qar\u015f\u0131s\u0131n\u0131 is\u0259 reflection il\u0259 rahatl\u0131qla ala bil\u0259rik. Lakin bunu bizim \u00fc\u00e7\u00fcn \u0259vv\u0259lc\u0259d\u0259n edib, AutoMapper package hal\u0131na sal\u0131blar.\",\n  \"\",\n  \"AutoMapper v\u0259 AutoMapper.Extensions.Microsoft.DependencyInjection package y\u00fckl\u0259nl\u0259nm\u0259lidir. Daha sonra h\u0259r Entity ad\u0131na uy\u011fun olaraq AutoMapper class yarad\u0131lmal\u0131, Package vasit\u0259sil\u0259 g\u0259l\u0259n Profile class\" + RTL_APOSTROPHE + \"dan miras almal\u0131d\u0131r. Bo\u015f constructor override edilm\u0259li, CreateMap<T,R>(); generic method i\u015f\u0259 sal\u0131nmal\u0131d\u0131r. Bu halda g\u00f6nd\u0259ril\u0259n T type(class) v\u0259 R type(class) daxilind\u0259ki property\" + RTL_APOSTROPHE + \"l\u0259r v\u0259 onlar\u0131n value\" + RTL_APOSTROPHE + \"lar\u0131 Assembly t\u0259r\u0259find\u0259n oxunur, T type daxilind\u0259ki property\" + RTL_APOSTROPHE + \"l\u0259rd\u0259n adlar\u0131 R type daxilind\u0259ki property\" + RTL_APOSTROPHE + \"l\u0259r il\u0259 \u00fcst-\u00fcst\u0259 d\u00fc\u015f\u0259nl\u0259r yeni R type object\" + RTL_APOSTROPHE + \"\u0259 m\u0259nims\u0259dilir v\u0259 geriy\u0259 qaytar\u0131l\u0131r. Burada adlar\u0131 eyni lakin datatype\" + RTL_APOSTROPHE + \"lar\u0131 f\u0259rqli olan property\" + RTL_APOSTROPHE + \"l\u0259r d\u0259 ola bil\u0259r. Bu halda runtime exception ba\u015f ver\u0259c\u0259k.\",\n  \"\",\n  \"ReverseMap(); method\" + RTL_APOSTROPHE + \"dan da istifad\u0259 ed\u0259r\u0259k h\u0259m Entity g\u00f6nd\u0259rib ViewModel & DTO object\" + RTL_APOSTROPHE + \"i h\u0259m d\u0259 ViewModel & DTO g\u00f6nd\u0259rib Entity objecti q\u0259bul ed\u0259 bil\u0259c\u0259yimizi bildiririk. \"\n];\n\n// Rewrite the intro paragraph's text in place (target[0]).\nparaIntro.insertText(newTexts[0], Word.InsertLocation.replace);\nawait context.sync();\n\n// paraBlank1 is already empty text (target[1]) \u2014 leave as is.\n\n// Insert target[2..6] (5 paragraphs) right after the blank paragraph,\n// preserving document order, then the existing paraBody paragraph is\n// rewritten to become target[7].\nlet anchor = paraBlank1;\nfor (let i = 1; i < newTexts.length; i++) {\n  anchor = anchor.insertParagraph(newTexts[i], Word.InsertLocation.after);\n  await context.sync();\n}\n\n// Finally, rewrite the old \"AutoMapper v\u0259 ...\" paragraph to hold only the\n// trailing \"Sonda services ...\" sentence (target[7]).\nparaBody.insertText(\n  \"Sonda services olaraq adland\u0131rd\u0131\u011f\u0131m\u0131z hiss\u0259d\u0259 g\u0259lib builder.Services.AddAutoMapper(typeof(mapperAssembly)); vasit\u0259sil\u0259 Mapper class\" + RTL_APOSTROPHE + \"lar\u0131n yerl\u0259\u015fdiyi Assembly\" + RTL_APOSTROPHE + \"ni tan\u0131daraq prosesi yekunla\u015fd\u0131r\u0131r\u0131q. \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document ($d below).\n#\n# Implements three textual revisions to the document (matching the\n# supplied OOXML diff):\n#   1. \"DbContext\" paragraph   - re-typing of \"class'dan\" (text unchanged).\n#   2. \"Configuration\" paragraph - re-typing of \"class'dan\" (text unchanged).\n#   3. \"AutoMapper\" section - expanded / rewritten explanation split\n#      across several paragraphs (new paragraphs inserted).\n\n$d = $word.ActiveDocument\n$apos = [char]0x2019   # ' RIGHT SINGLE QUOTATION MARK, used throughout the doc\n\nfunction Get-ParagraphIndexContaining {\n    param([string]$needle, [int]$startAt = 1)\n    $count = $d.Paragraphs.Count\n    for ($i = $startAt; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text.Contains($needle)) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Retype-Substring {\n    # Re-enters $needle inside the paragraph whose text contains\n    # $paraNeedle, leaving the visible text unchanged.\n    param([string]$paraNeedle, [string]$needle)\n    $idx = Get-ParagraphIndexContaining $paraNeedle\n    if ($idx -lt 0) {\n        throw \"Paragraph containing '$paraNeedle' not found\"\n    }\n    $p = $d.Paragraphs.Item($idx)\n    $rng = $p.Range.Duplicate\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $needle\n    $find.MatchCase = $true\n    $find.MatchWildcards = $false\n    $find.Execute() | Out-Null\n    if ($find.Found) {\n        $rng.Text = $needle\n    } else {\n        throw \"Substring '$needle' not found in paragraph '$paraNeedle'\"\n    }\n}\n\n# 1) \"\u0130lk olaraq DbContext class'dan miras alan ...\" paragraph.\nRetype-Substring \"\u0130lk olaraq DbContext\" (\"class\" + $apos + \"dan\")\n\n# 2) \".cs file daxilind\u0259 ... Configuration class'dan istifad\u0259 etm\u0259kdir.\" paragraph.\nRetype-Substring \"Configuration\" (\"class\" + $apos + \"dan\")\n\n# ---------------------------------------------------------------------\n# 3) AutoMapper section rewrite.\n# ---------------------------------------------------------------------\n$idxIntro = Get-ParagraphIndexContaining \"Biz project daxilind\u0259 Front\"\nif ($idxIntro -lt 0) { throw \"Could not locate AutoMapper intro paragraph\" }\n$idxBody = Get-ParagraphIndexContaining \"AutoMapper v\u0259 AutoMapper.Extensions.Microsoft.DependencyInjection\" $idxIntro\nif ($idxBody -lt 0) { throw \"Could not locate AutoMapper body paragraph\" }\n\n$text0 = \"Biz security c\u0259h\u0259td\u0259n project daxilind\u0259 Front\" + $apos + \"a database\" + $apos + \"d\u0259n g\u0259l\u0259n Entity object\" + $apos + \"l\u0259ri bir ba\u015fa olaraq g\u00f6nd\u0259rm\u0259yimiz t\u0259hl\u00fck\u0259lidir. \u00c7\u00fcnki biz bu object\" + $apos + \"l\u0259ri bir ba\u015fa \u00f6t\u00fcrs\u0259k \u0259l\u00e7atan olmamal\u0131 bir s\u0131ra datalar\u0131 da (property\" + $apos + \"l\u0259ri) g\u00f6nd\u0259rmi\u015f olar\u0131q. Bunun qar\u015f\u0131s\u0131n\u0131 almaq \u00fc\u00e7\u00fcn ViewModel v\u0259 ya DTO class\" + $apos + \"lar\u0131 yarad\u0131b m\u00fc\u0259yy\u0259n olunmu\u015f property\" + $apos + \"l\u0259ri saxlay\u0131r\u0131q.\"\n\n$text2 = \"H\u0259r d\u0259f\u0259 datan\u0131n \u00f6t\u00fcr\u00fclm\u0259si v\u0259 ya q\u0259bul olunmas\u0131 zaman\u0131 yeni ViewModel & DTO v\u0259 ya Entity object yarad\u0131lmas\u0131 v\u0259 d\u0259y\u0259rl\u0259rin m\u0259nims\u0259dilm\u0259si h\u0259m kod t\u0259krar\u0131na g\u0259tirib \u00e7\u0131xar\u0131r h\u0259m d\u0259 vaxt itkisin\u0259. Bu problemin qar\u015f\u0131s\u0131n\u0131 is\u0259 reflection il\u0259 rahatl\u0131qla ala bil\u0259rik. Lakin bunu bizim \u00fc\u00e7\u00fcn \u0259vv\u0259lc\u0259d\u0259n edib, AutoMapper package hal\u0131na sal\u0131blar.\"\n\n$text4 = \"AutoMapper v\u0259 AutoMapper.Extensions.Microsoft.DependencyInjection package y\u00fckl\u0259nl\u0259nm\u0259lidir. Daha sonra h\u0259r Entity ad\u0131na uy\u011fun olaraq AutoMapper class yarad\u0131lmal\u0131, Package vasit\u0259sil\u0259 g\u0259l\u0259n Profile class\" + $apos + \"dan miras almal\u0131d\u0131r. Bo\u015f constructor override edilm\u0259li, CreateMap<T,R>(); generic method i\u015f\u0259 sal\u0131nmal\u0131d\u0131r. Bu halda g\u00f6nd\u0259ril\u0259n T type(class) v\u0259 R type(class) daxilind\u0259ki property\" + $apos + \"l\u0259r v\u0259 onlar\u0131n value\" + $apos + \"lar\u0131 Assembly t\u0259r\u0259find\u0259n oxunur, T type daxilind\u0259ki property\" + $apos + \"l\u0259rd\u0259n adlar\u0131 R type daxilind\u0259ki property\" + $apos + \"l\u0259r il\u0259 \u00fcst-\u00fcst\u0259 d\u00fc\u015f\u0259nl\u0259r yeni R type object\" + $apos + \"\u0259 m\u0259nims\u0259dilir v\u0259 geriy\u0259 qaytar\u0131l\u0131r. Burada adlar\u0131 eyni lakin datatype\" + $apos + \"lar\u0131 f\u0259rqli olan property\" + $apos + \"l\u0259r d\u0259 ola bil\u0259r. Bu halda runtime exception ba\u015f ver\u0259c\u0259k.\"\n\n$text6 = \"ReverseMap(); method\" + $apos + \"dan da istifad\u0259 ed\u0259r\u0259k h\u0259m Entity g\u00f6nd\u0259rib ViewModel & DTO object\" + $apos + \"i h\u0259m d\u0259 ViewModel & DTO g\u00f6nd\u0259rib Entity objecti q\u0259bul ed\u0259 bil\u0259c\u0259yimizi bildiririk. \"\n\n$text7 = \"Sonda services olaraq adland\u0131rd\u0131\u011f\u0131m\u0131z hiss\u0259d\u0259 g\u0259lib builder.Services.AddAutoMapper(typeof(mapperAssembly)); vasit\u0259sil\u0259 Mapper class\" + $apos + \"lar\u0131n yerl\u0259\u015fdiyi Assembly\" + $apos + \"ni tan\u0131daraq prosesi yekunla\u015fd\u0131r\u0131r\u0131q. \"\n\n# Rewrite the intro paragraph text in place (-> target[0]).\n$paraIntro = $d.Paragraphs.Item($idxIntro)\n$paraIntro.Range.Text = $text0\n\n# The paragraph right after the intro is already empty (-> target[1]); leave it.\n$idxBlank1 = $idxIntro + 1\n$paraBlank1 = $d.Paragraphs.Item($idxBlank1)\n\n# Insert target[2..6] (5 new paragraphs) right after the blank paragraph.\n$newTexts = @($text2, \"\", $text4, \"\", $text6)\n\n$anchor = $paraBlank1\nforeach ($t in $newTexts) {\n    $anchor.Range.InsertParagraphAfter()\n    $idxBlank1 = $idxBlank1 + 1\n    $anchor = $d.Paragraphs.Item($idxBlank1)\n    if ($t -ne \"\") {\n        $anchor.Range.Text = $t\n    }\n}\n\n# Finally, rewrite the old \"AutoMapper v\u0259 ...\" paragraph (shifted by the 5\n# freshly-inserted paragraphs) to hold only the trailing \"Sonda services ...\"\n# sentence (-> target[7]).\n$idxBody = $idxBody + 5\n$paraBody = $d.Paragraphs.Item($idxBody)\n$paraBody.Range.Text = $text7\n\nWrite-Output \"done\"\n"}
